$wb = $excel.ActiveWorkbook

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N")

# Adds a new scrim-log row by copying an existing template row (to inherit
# fonts/fills/borders) and then overwriting the cell values with the new
# row's data. $equipoStyleRow, if non-zero, is a row whose G-column style
# (bold + correct fill for "Equipo 1" vs "Equipo 2") is copied onto the new
# row's G cell when it doesn't match the template row's G style.
function Add-ScrimRow($ws, $templateRow, $newRow, $values, $equipoStyleRow) {
    $srcRange = $ws.Range("A" + $templateRow + ":N" + $templateRow)
    $dstRange = $ws.Range("A" + $newRow + ":N" + $newRow)
    $srcRange.Copy($dstRange) | Out-Null

    if ($equipoStyleRow -ne 0) {
        $ws.Range("G" + $equipoStyleRow).Copy($ws.Range("G" + $newRow)) | Out-Null
    }

    foreach ($col in $cols) {
        $ws.Range($col + $newRow).Value = $values[$col]
    }
}

# ---------------------------------------------------------------------
# Sheet "Hot Potato": append rows 81-82 (both "Equipo 2", same style as
# template row 80, so no G restyle needed).
# ---------------------------------------------------------------------
$wsHotPotato = $wb.Worksheets.Item("Hot Potato")

$row81 = @{
    "A"="CHARLIE"; "B"="NITA"; "C"="LILY"; "D"="EMZ"; "E"="KAZE"; "F"="MICO";
    "G"="Equipo 2"; "H"="NHG|Xemp"; "I"="KCP|Tyrant"; "J"="KCP|Fade";
    "K"="SSG|bobby"; "L"="SSG|chino"; "M"="CODE|OG"; "N"="20250724T223516.000Z"
}
Add-ScrimRow $wsHotPotato 80 81 $row81 0

$row82 = @{
    "A"="CHARLIE"; "B"="NITA"; "C"="LILY"; "D"="EMZ"; "E"="KAZE"; "F"="MICO";
    "G"="Equipo 2"; "H"="NHG|Xemp"; "I"="KCP|Tyrant"; "J"="KCP|Fade";
    "K"="SSG|bobby"; "L"="SSG|chino"; "M"="CODE|OG"; "N"="20250724T223402.000Z"
}
Add-ScrimRow $wsHotPotato 81 82 $row82 0

# ---------------------------------------------------------------------
# Sheet "Open Business": append rows 72-78. Template row 71 is styled
# "Equipo 2" (G s=5); row 69 is styled "Equipo 1" (G s=4) and is used as
# the style donor whenever a new row needs "Equipo 1".
# ---------------------------------------------------------------------
$wsOpenBusiness = $wb.Worksheets.Item("Open Business")

$row72 = @{
    "A"="AMBER"; "B"="MOE"; "C"="CORDELIUS"; "D"="KIT"; "E"="ASH"; "F"="FINX";
    "G"="Equipo 1"; "H"="TE|Rafikii"; "I"="TE|Ezlivi"; "J"="TE|Belal";
    "K"="TRB|Zeus 解開"; "L"="TRB|R B M"; "M"="TRB|Lxffy"; "N"="20250724T223419.000Z"
}
Add-ScrimRow $wsOpenBusiness 71 72 $row72 69

$row73 = @{
    "A"="AMBER"; "B"="MOE"; "C"="CORDELIUS"; "D"="KIT"; "E"="ASH"; "F"="FINX";
    "G"="Equipo 1"; "H"="TE|Rafikii"; "I"="TE|Ezlivi"; "J"="TE|Belal";
    "K"="TRB|Zeus 解開"; "L"="TRB|R B M"; "M"="TRB|Lxffy"; "N"="20250724T223216.000Z"
}
Add-ScrimRow $wsOpenBusiness 72 73 $row73 0

$row74 = @{
    "A"="KAZE"; "B"="LUMI"; "C"="GRIFF"; "D"="GRAY"; "E"="CORDELIUS"; "F"="DOUG";
    "G"="Equipo 1"; "H"="NHG|Xemp"; "I"="KCP|Tyrant"; "J"="KCP|Fade";
    "K"="SSG|bobby"; "L"="SSG|chino"; "M"="CODE|OG"; "N"="20250724T222525.000Z"
}
Add-ScrimRow $wsOpenBusiness 73 74 $row74 69

$row75 = @{
    "A"="KAZE"; "B"="LUMI"; "C"="GRIFF"; "D"="GRAY"; "E"="CORDELIUS"; "F"="DOUG";
    "G"="Equipo 2"; "H"="NHG|Xemp"; "I"="KCP|Tyrant"; "J"="KCP|Fade";
    "K"="SSG|bobby"; "L"="SSG|chino"; "M"="CODE|OG"; "N"="20250724T222340.000Z"
}
Add-ScrimRow $wsOpenBusiness 74 75 $row75 71

$row76 = @{
    "A"="KAZE"; "B"="LUMI"; "C"="GRIFF"; "D"="GRAY"; "E"="CORDELIUS"; "F"="DOUG";
    "G"="Equipo 1"; "H"="NHG|Xemp"; "I"="KCP|Tyrant"; "J"="KCP|Fade";
    "K"="SSG|bobby"; "L"="SSG|chino"; "M"="CODE|OG"; "N"="20250724T222120.000Z"
}
Add-ScrimRow $wsOpenBusiness 75 76 $row76 69

$row77 = @{
    "A"="HANK"; "B"="KAZE"; "C"="BARLEY"; "D"="LOU"; "E"="LARRY & LAWRIE"; "F"="DRACO";
    "G"="Equipo 2"; "H"="NHG|Xemp"; "I"="KCP|Fade"; "J"="KCP|Tyrant";
    "K"="SSG|bobby"; "L"="CODE|OG"; "M"="SSG|chino"; "N"="20250724T221548.000Z"
}
Add-ScrimRow $wsOpenBusiness 76 77 $row77 71

$row78 = @{
    "A"="HANK"; "B"="KAZE"; "C"="BARLEY"; "D"="LOU"; "E"="LARRY & LAWRIE"; "F"="DRACO";
    "G"="Equipo 2"; "H"="NHG|Xemp"; "I"="KCP|Fade"; "J"="KCP|Tyrant";
    "K"="SSG|bobby"; "L"="CODE|OG"; "M"="SSG|chino"; "N"="20250724T221304.000Z"
}
Add-ScrimRow $wsOpenBusiness 77 78 $row78 0
